# Updated dev tasks for use case 3 - Near store activity
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Row 2 (S.No. 1)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 'Create addon "storeActivityAddOn" and install it on bncwebservices'
$ws.Range("C2").Value = "Swarnima/Swapnil"
$ws.Range("D2").Value = "23/03"
$ws.Range("E2").Value = "23/03"

# ---------------------------------------------------------------------------
# Row 3 (S.No. 2)
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "Create an item ""storeActivity"" and add following attributes to it.
1. storeId
2. customerId
3. storeVisitDate
4. storeEntryTime
5.storeExitTime
6. timeSpentInStore (save time in minutes, to be calculated only when storeEntryTime and storeExitTime are available)"
$ws.Range("C3").Value = "Swapnil"
$ws.Range("D3").Value = "24/03"
$ws.Range("E3").Value = "24/03"

# ---------------------------------------------------------------------------
# Row 4 (S.No. 3)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Create a web service that takes following prameters in JSON format and save the data in hybris.
1. storeId
2. customerId
3. storeVisitDate
4. storeEntryTime
5.storeExitTime
6. timeSpentInStore (save time in minutes, to be calculated only when storeEntryTime and storeExitTime are available)"
$ws.Range("C4").Value = "Swapnil"
$ws.Range("D4").Value = "24/03"
$ws.Range("E4").Value = "25/03"

# ---------------------------------------------------------------------------
# Row 5 (S.No. 4)
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 'Create service and dao classes to get most visited stores by customer. The "storeActivity" item (created in step 2) is to be queried to get most visited stores'
$ws.Range("C5").Value = "Swarnima"
$ws.Range("D5").Value = "24/03"
$ws.Range("E5").Value = "25/03"

# ---------------------------------------------------------------------------
# Row 6 (S.No. 5)
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = 'Create service and dao classes to get loyal customers. The "storeActivity" item is to be queried to get the customers who visit stores the most'
$ws.Range("C6").Value = "Swarnima"
$ws.Range("D6").Value = "26/03"
$ws.Range("E6").Value = "27/03"

# ---------------------------------------------------------------------------
# Row 7 (S.No. 6) - new row
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("B7").WrapText = $true
$ws.Range("B7").Value = 'Create an item type "beacon" that contains following String type attributes.
1. beaconId
2. majorId
3. minorId'
$ws.Range("C7").Value = "Swapnil"
$ws.Range("D7").Value = "26/03"
$ws.Range("E7").Value = "26/03"

# ---------------------------------------------------------------------------
# Row 8 (S.No. 7) - new row
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 7
$ws.Range("B8").WrapText = $true
$ws.Range("B8").Value = "Customize product item type by adding following attribute.
1. popularityCount (int type)"
$ws.Range("C8").Value = "Swapnil"
$ws.Range("D8").Value = "26/03"
$ws.Range("E8").Value = "26/03"

# ---------------------------------------------------------------------------
# Row 9 (S.No. 8) - new row
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = 8
$ws.Range("B9").WrapText = $true
$ws.Range("B9").Value = "Create a one to many relation between beacon and products. And add some sample data in it (create impex)"
$ws.Range("C9").Value = "Swapnil"
$ws.Range("D9").Value = "26/03"
$ws.Range("E9").Value = "26/03"

# ---------------------------------------------------------------------------
# Row 10 (S.No. 9) - new row
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").WrapText = $true
$ws.Range("B10").Value = "Create a webservice to update the popularityCount of the product. The web service will accept the beaconId, majorId, minorId (in JSON format). Will use the request data to get the product from the relation defined in point 8 above and increase the popularity count for the product."
$ws.Range("C10").Value = "Swapnil"
$ws.Range("D10").Value = "26/03"
$ws.Range("E10").Value = "31/03"

# ---------------------------------------------------------------------------
# Row 11 (S.No. 10) - new row (real dates)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 10
$ws.Range("B11").WrapText = $true
$ws.Range("B11").Value = "Create service and dao classes to get top 5 popular products based on the popularity count. The count ""5"" should be made configurable."
$ws.Range("C11").Value = "Swapnil"
$ws.Range("D11:E11").NumberFormat = "d-mmm"
$ws.Range("D11").Value = (Get-Date -Year 2015 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E11").Value = (Get-Date -Year 2015 -Month 4 -Day 2 -Hour 0 -Minute 0 -Second 0)

# ---------------------------------------------------------------------------
# Row 12 (S.No. 11) - new row (real dates)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 11
$ws.Range("B12").WrapText = $true
$ws.Range("B12").Value = 'create an item "storeCustomer" and add following attributes to it.
1. weight
2. height
3. Age
4. Gender
'
$ws.Range("C12").Value = "Swapnil"
$ws.Range("D12:E12").NumberFormat = "d-mmm"
$ws.Range("D12").Value = (Get-Date -Year 2015 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E12").Value = (Get-Date -Year 2015 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0)

# ---------------------------------------------------------------------------
# Row 13 (S.No. 12) - new row (real dates)
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = 12
$ws.Range("B13").WrapText = $true
$ws.Range("B13").Value = "Create impex to store some sample data in storeCustomer point 11 above."
$ws.Range("C13").Value = "Swapnil"
$ws.Range("D13:E13").NumberFormat = "d-mmm"
$ws.Range("D13").Value = (Get-Date -Year 2015 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E13").Value = (Get-Date -Year 2015 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0)

# ---------------------------------------------------------------------------
# Row 14 (S.No. 13) - new row
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = 13
$ws.Range("B14").WrapText = $true
$ws.Range("B14").Value = "Create a WCMS page for the Activity dashboard, impex creation."
$ws.Range("C14").Value = "Swarnima"
$ws.Range("D14").Value = "30/03"
$ws.Range("E14").Value = "30/03"

# ---------------------------------------------------------------------------
# Row 15 (S.No. 14) - new row (E is a real date)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = 14
$ws.Range("B15").WrapText = $true
$ws.Range("B15").Value = "Create controller, facades, service classes to fetch model data for store customer profile section (point 11), popular products (based on product popularity count, point 7), most visited stores (point 4), loyal customers (from point 5), spent time (query storeActivity item, point 2)"
$ws.Range("C15").Value = "Swarnima"
$ws.Range("D15").Value = "31/03"
$ws.Range("E15").NumberFormat = "d-mmm"
$ws.Range("E15").Value = (Get-Date -Year 2015 -Month 4 -Day 2 -Hour 0 -Minute 0 -Second 0)

# ---------------------------------------------------------------------------
# Row 16 (S.No. 15) - new row (real dates)
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 15
$ws.Range("B16").WrapText = $true
$ws.Range("B16").Value = 'To get weather information, there are two approaches. 
1) The weather data is passed by IOS app, if so we can use it to dispplay on our page.
2) To use java APIs to fetch weather data based on city or zip code. Please see the sample code http://code.aksingh.net/owm-japis/src'
$ws.Range("C16").Value = "Swarnima"
$ws.Range("D16:E16").NumberFormat = "d-mmm"
$ws.Range("D16").Value = (Get-Date -Year 2015 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E16").Value = (Get-Date -Year 2015 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0)

# ---------------------------------------------------------------------------
# Apply an AutoFilter on the Owner column (C1:C16)
# ---------------------------------------------------------------------------
$ws.Range("C1:C16").AutoFilter(1)

# Update the active selection to match the final edited cell, like the author did
$ws.Range("F15").Select()
